$wb = $excel.ActiveWorkbook

# --- "tip deflection results" sheet: insert new row for test #10 ---
$wsTip = $wb.Worksheets.Item("tip deflection results")

# Insert a new row above current row 6 (shifts old rows 6-12 down to 7-13)
$wsTip.Rows.Item(6).Insert()

$wsTip.Cells.Item(6, 1).Value = 10
$wsTip.Cells.Item(6, 2).Value = 0.274
$wsTip.Cells.Item(6, 3).Value = 0.75
$wsTip.Cells.Item(6, 4).Value = 0.70512802088869997

$wsTip.Range("D6").Select()

# --- "test matrix" sheet: update status cells ---
$wsMatrix = $wb.Worksheets.Item("test matrix")

$wsMatrix.Cells.Item(12, 1).Value = "completed"
$wsMatrix.Cells.Item(19, 1).Value = "in progress"

$wsMatrix.Range("A20").Select()
